# Weekly update: a new price observation was inserted for this market/product
# series at row 35 ("Fecha" 44484, i.e. 2021-10-15), which pushes every
# subsequent observation (previously rows 35-80) down by one row (now 36-81).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 35; Excel shifts rows 35:80 down to 36:81
# and the new blank row inherits the formatting (incl. the date style on
# column D) from the row above it.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44484
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108002
$ws.Range("J35").Value = "Mango"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 208
$ws.Range("N35").Value = 7000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 7481
$ws.Range("Q35").Value = "$/bandeja 4 kilos"
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1870
$ws.Range("T35").Value = 4
